$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set plain text/string values (safe from numeric auto-detection) ---
$ws.Cells.Item(2, 4).Value = '64.250.76'
$ws.Cells.Item(2, 5).Value = '  -2.48%  '
$ws.Cells.Item(3, 4).Value = '3.180.69'
$ws.Cells.Item(3, 5).Value = '  -3.14%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 5).Value = '  -2.46%  '
$ws.Cells.Item(6, 5).Value = '  -6.11%  '
$ws.Cells.Item(7, 5).Value = '  -5.13%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 4).Value = '3.181.71'
$ws.Cells.Item(9, 5).Value = '  -2.88%  '
$ws.Cells.Item(10, 5).Value = '  -3.24%  '
$ws.Cells.Item(11, 5).Value = '  -1.19%  '
$ws.Cells.Item(12, 5).Value = '  -4.19%  '
$ws.Cells.Item(13, 4).Value = '3.724.00'
$ws.Cells.Item(13, 5).Value = '  -3.36%  '
$ws.Cells.Item(14, 5).Value = '  -1.85%  '
$ws.Cells.Item(15, 4).Value = '64.334.19'
$ws.Cells.Item(15, 5).Value = '  -2.45%  '
$ws.Cells.Item(17, 5).Value = '  -2.38%  '
$ws.Cells.Item(18, 4).Value = '3.171.11'
$ws.Cells.Item(18, 5).Value = '  -5.44%  '
$ws.Cells.Item(19, 5).Value = '  -1.84%  '
$ws.Cells.Item(20, 5).Value = '  -2.88%  '
$ws.Cells.Item(21, 5).Value = '  -2.47%  '
$ws.Cells.Item(22, 5).Value = '  -3.88%  '
$ws.Cells.Item(23, 5).Value = '  +0.03%  '
$ws.Cells.Item(24, 5).Value = '  -2.10%  '
$ws.Cells.Item(25, 5).Value = '  +2.59%  '
$ws.Cells.Item(26, 5).Value = '  -3.77%  '
$ws.Cells.Item(27, 5).Value = '  -5.38%  '
$ws.Cells.Item(28, 5).Value = '  -1.38%  '
$ws.Cells.Item(29, 5).Value = '  -0.22%  '
$ws.Cells.Item(30, 5).Value = '  -5.69%  '
$ws.Cells.Item(31, 5).Value = '  -1.74%  '
$ws.Cells.Item(32, 5).Value = '  -0.07%  '
$ws.Cells.Item(33, 5).Value = '  -2.75%  '
$ws.Cells.Item(34, 5).Value = '  -3.52%  '
$ws.Cells.Item(35, 5).Value = '  -4.28%  '
$ws.Cells.Item(36, 5).Value = '  -2.27%  '
$ws.Cells.Item(37, 5).Value = '  -4.12%  '
$ws.Cells.Item(38, 4).Value = '2.704.04'
$ws.Cells.Item(38, 5).Value = '  -3.85%  '
$ws.Cells.Item(39, 5).Value = '  -6.04%  '
$ws.Cells.Item(40, 5).Value = '  -7.45%  '
$ws.Cells.Item(41, 5).Value = '  -3.21%  '
$ws.Cells.Item(42, 5).Value = '  -2.91%  '
$ws.Cells.Item(43, 5).Value = '  -7.04%  '
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 5).Value = '  -4.24%  '
$ws.Cells.Item(45, 2).Value = 'Hedera'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(45, 5).Value = '  -5.30%  '
$ws.Cells.Item(46, 5).Value = '  -4.97%  '
$ws.Cells.Item(47, 5).Value = '  -5.41%  '
$ws.Cells.Item(48, 5).Value = '  -1.96%  '
$ws.Cells.Item(49, 5).Value = '  -10.11%  '
$ws.Cells.Item(50, 5).Value = '  -4.30%  '
$ws.Cells.Item(51, 5).Value = '  +0.01%  '

# --- Values that look numeric must be forced to Text so formatting (e.g. trailing zeros) is preserved ---
$numericLookingCells = @(
    $ws.Cells.Item(5, 4),
    $ws.Cells.Item(6, 4),
    $ws.Cells.Item(7, 4),
    $ws.Cells.Item(10, 4),
    $ws.Cells.Item(11, 4),
    $ws.Cells.Item(16, 4),
    $ws.Cells.Item(17, 4),
    $ws.Cells.Item(19, 4),
    $ws.Cells.Item(20, 4),
    $ws.Cells.Item(22, 4),
    $ws.Cells.Item(23, 4),
    $ws.Cells.Item(24, 4),
    $ws.Cells.Item(25, 4),
    $ws.Cells.Item(28, 4),
    $ws.Cells.Item(29, 4),
    $ws.Cells.Item(30, 4),
    $ws.Cells.Item(31, 4),
    $ws.Cells.Item(33, 4),
    $ws.Cells.Item(34, 4),
    $ws.Cells.Item(35, 4),
    $ws.Cells.Item(36, 4),
    $ws.Cells.Item(40, 4),
    $ws.Cells.Item(41, 4),
    $ws.Cells.Item(42, 4),
    $ws.Cells.Item(43, 4),
    $ws.Cells.Item(44, 4),
    $ws.Cells.Item(45, 4),
    $ws.Cells.Item(46, 4),
    $ws.Cells.Item(47, 4),
    $ws.Cells.Item(49, 4),
    $ws.Cells.Item(50, 4),
    $ws.Cells.Item(51, 4)
)
foreach ($cell in $numericLookingCells) {
    $cell.NumberFormat = "@"
}

$ws.Cells.Item(5, 4).Value = '569.86'
$ws.Cells.Item(6, 4).Value = '169.18'
$ws.Cells.Item(7, 4).Value = '0.609'
$ws.Cells.Item(10, 4).Value = '0.121'
$ws.Cells.Item(11, 4).Value = '6.70'
$ws.Cells.Item(16, 4).Value = '25.40'
$ws.Cells.Item(17, 4).Value = '0.0000159'
$ws.Cells.Item(19, 4).Value = '418.89'
$ws.Cells.Item(20, 4).Value = '12.82'
$ws.Cells.Item(22, 4).Value = '7.08'
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(24, 4).Value = '70.11'
$ws.Cells.Item(25, 4).Value = '0.203'
$ws.Cells.Item(28, 4).Value = '8.80'
$ws.Cells.Item(29, 4).Value = '0.996'
$ws.Cells.Item(30, 4).Value = '1.84'
$ws.Cells.Item(31, 4).Value = '21.84'
$ws.Cells.Item(33, 4).Value = '5.02'
$ws.Cells.Item(34, 4).Value = '6.35'
$ws.Cells.Item(35, 4).Value = '1.14'
$ws.Cells.Item(36, 4).Value = '155.57'
$ws.Cells.Item(40, 4).Value = '24.60'
$ws.Cells.Item(41, 4).Value = '4.18'
$ws.Cells.Item(42, 4).Value = '38.86'
$ws.Cells.Item(43, 4).Value = '0.710'
$ws.Cells.Item(44, 4).Value = '5.69'
$ws.Cells.Item(45, 4).Value = '0.0624'
$ws.Cells.Item(46, 4).Value = '21.93'
$ws.Cells.Item(47, 4).Value = '297.05'
$ws.Cells.Item(49, 4).Value = '2.06'
$ws.Cells.Item(50, 4).Value = '0.0994'
$ws.Cells.Item(51, 4).Value = '0.999'

# Reset number format back to the sheet default so the cell style matches the original (unstyled) cells
foreach ($cell in $numericLookingCells) {
    $cell.Style = "Normal"
}